$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2961.9048
$ws.Range("I64").Value = 2938.4614
$ws.Range("K64").Value = 2938.4614
$ws.Range("M64").Value = -2690.4614
$ws.Range("H67").Value = 2961.9048
$ws.Range("I67").Value = 2938.4614
$ws.Range("K67").Value = 2938.4614
$ws.Range("M67").Value = -2080.4614
$ws.Range("H76").Value = 3261.9048
$ws.Range("I76").Value = 3270.5881
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 3270.5881
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -2955.5881
$ws.Range("N76").Value = -3855
$ws.Range("H79").Value = 3261.9048
$ws.Range("I79").Value = 3270.5881
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 3270.5881
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -2178.5881
$ws.Range("N79").Value = -5409
$ws.Range("H106").Value = 22225068
$ws.Range("I106").Value = 30002968
$ws.Range("K106").Value = 30002968
$ws.Range("M106").Value = -30002337
$ws.Range("H137").Value = 2173.3044
$ws.Range("I137").Value = 2291.3635
$ws.Range("J137").Value = 2065.0833
$ws.Range("K137").Value = 6874.0905
$ws.Range("L137").Value = 6195.249899999999
$ws.Range("M137").Value = -4324.0905
$ws.Range("N137").Value = -11295.2499
$ws.Range("H138").Value = 3334.4387
$ws.Range("I138").Value = 4050.0908
$ws.Range("J138").Value = 3243.954
$ws.Range("K138").Value = 12150.2724
$ws.Range("L138").Value = 9731.862000000001
$ws.Range("M138").Value = -7010.2724
$ws.Range("N138").Value = -20011.862
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 40327.72
$ws.Range("I63").Value = 183093.6
$ws.Range("J63").Value = 4636.25
$ws.Range("K63").Value = 183093.6
$ws.Range("L63").Value = 4636.25
$ws.Range("M63").Value = -182407.6
$ws.Range("N63").Value = -6008.25
$ws.Range("H66").Value = 40327.72
$ws.Range("I66").Value = 183093.6
$ws.Range("J66").Value = 4636.25
$ws.Range("K66").Value = 915468
$ws.Range("L66").Value = 23181.25
$ws.Range("M66").Value = -912036
$ws.Range("N66").Value = -30045.25
$ws.Range("H88").Value = 2701.5
$ws.Range("I88").Value = 2602
$ws.Range("K88").Value = 2602
$ws.Range("M88").Value = -2196
$ws.Range("H91").Value = 2701.5
$ws.Range("I91").Value = 2602
$ws.Range("K91").Value = 2602
$ws.Range("M91").Value = -1198
$ws.Range("H93").Value = 65349.715
$ws.Range("J93").Value = 65349.715
$ws.Range("L93").Value = 65349.715
$ws.Range("N93").Value = -70341.715
$ws.Range("H97").Value = 1258.1
$ws.Range("I97").Value = 1063.3334
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 1063.3334
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -567.3334
$ws.Range("N97").Value = -4003
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2536.7144
$ws.Range("I86").Value = 2226.7856
$ws.Range("J86").Value = 3156.5715
$ws.Range("K86").Value = 2226.7856
$ws.Range("L86").Value = 3156.5715
$ws.Range("M86").Value = -1103.7856
$ws.Range("N86").Value = -5402.5715
$ws.Range("H89").Value = 2536.7144
$ws.Range("I89").Value = 2226.7856
$ws.Range("J89").Value = 3156.5715
$ws.Range("K89").Value = 11133.928
$ws.Range("L89").Value = 15782.8575
$ws.Range("M89").Value = -5517.928
$ws.Range("N89").Value = -27014.8575
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 17858982
$ws.Range("I105").Value = 20835184
$ws.Range("J105").Value = 1780
$ws.Range("K105").Value = 20835184
$ws.Range("L105").Value = 1780
$ws.Range("M105").Value = -20833437
$ws.Range("N105").Value = -5274
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 659.8077
$ws.Range("I5").Value = 525.36
$ws.Range("J5").Value = 899.8929000000001
$ws.Range("K5").Value = 1576.08
$ws.Range("L5").Value = 2699.6787
$ws.Range("M5").Value = -1464.08
$ws.Range("N5").Value = -2923.6787
$ws.Range("H122").Value = 2552.0715
$ws.Range("J122").Value = 3714.2222
$ws.Range("L122").Value = 33427.99980000001
$ws.Range("N122").Value = -38327.99980000001
$ws.Range("H135").Value = 659.8077
$ws.Range("I135").Value = 525.36
$ws.Range("J135").Value = 899.8929000000001
$ws.Range("K135").Value = 4728.24
$ws.Range("L135").Value = 8099.0361
$ws.Range("M135").Value = -2193.24
$ws.Range("N135").Value = -13169.0361
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5955.875
$ws.Range("I43").Value = 1129.6
$ws.Range("J43").Value = 13999.667
$ws.Range("K43").Value = 1129.6
$ws.Range("L43").Value = 13999.667
$ws.Range("M43").Value = -978.5999999999999
$ws.Range("N43").Value = -14301.667
$ws.Range("H57").Value = 18999.77
$ws.Range("J57").Value = 18999.77
$ws.Range("L57").Value = 18999.77
$ws.Range("N57").Value = -20639.77
$ws.Range("H70").Value = 7688.242
$ws.Range("I70").Value = 8579.041999999999
$ws.Range("J70").Value = 5312.778
$ws.Range("K70").Value = 8579.041999999999
$ws.Range("L70").Value = 5312.778
$ws.Range("M70").Value = -8309.041999999999
$ws.Range("N70").Value = -5852.778
$ws.Range("H73").Value = 7688.242
$ws.Range("I73").Value = 8579.041999999999
$ws.Range("J73").Value = 5312.778
$ws.Range("K73").Value = 8579.041999999999
$ws.Range("L73").Value = 5312.778
$ws.Range("M73").Value = -7643.041999999999
$ws.Range("N73").Value = -7184.778
$ws.Range("H80").Value = 29345746
$ws.Range("I80").Value = 35134896
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 35134896
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -35133898
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 29345746
$ws.Range("I83").Value = 35134896
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 175674480
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -175669488
$ws.Range("N83").Value = -2009984
$ws.Range("H97").Value = 2710
$ws.Range("I97").Value = 2745.5557
$ws.Range("K97").Value = 2745.5557
$ws.Range("M97").Value = -2249.5557
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3784.6
$ws.Range("I81").Value = 3689.1667
$ws.Range("J81").Value = 4166.3335
$ws.Range("K81").Value = 7378.3334
$ws.Range("L81").Value = 8332.666999999999
$ws.Range("M81").Value = -6317.3334
$ws.Range("N81").Value = -10454.667
$ws.Range("H84").Value = 3784.6
$ws.Range("I84").Value = 3689.1667
$ws.Range("J84").Value = 4166.3335
$ws.Range("K84").Value = 36891.667
$ws.Range("L84").Value = 41663.335
$ws.Range("M84").Value = -31587.667
$ws.Range("N84").Value = -52271.335
$ws.Range("H96").Value = 4310
$ws.Range("I96").Value = 4131.6665
$ws.Range("J96").Value = 4666.6665
$ws.Range("K96").Value = 4131.6665
$ws.Range("L96").Value = 4666.6665
$ws.Range("M96").Value = -2758.6665
$ws.Range("N96").Value = -7412.6665
$ws.Range("H113").Value = 1589.5385
$ws.Range("I113").Value = 1646.4
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 4939.200000000001
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = -2769.200000000001
$ws.Range("N113").Value = -8540
$ws.Range("H133").Value = 40315
$ws.Range("J133").Value = 40315
$ws.Range("L133").Value = 40315
$ws.Range("N133").Value = -50435
